$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume snapshot refresh (15-1-2023 09:xx UTC run).
# Cells already hold plain-text numeric-looking values (no cell style), so
# force Text number format before writing to keep them as text, matching
# the source data which stores these as inline strings, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "294.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.53%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.78%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.088"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.68%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07368"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.44%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.681"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.44%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.763"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.68%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.662"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "8.90%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9243"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.45%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1672"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.45%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07120"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.52%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07864"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.02%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02994"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.53%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09893"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.09%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.85%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006282"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.50%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.455"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.226"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.75%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3279"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.17%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1350"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.69%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.571"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.45%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04644"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.76%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.68%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.85%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004419"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.86%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.18%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001880"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.49%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01660"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04392"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.53%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007075"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.08%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.10%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002103"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.16%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01101"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-23.13%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005984"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.19%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.930"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.96%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-16.36%"
